$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 4076;  "C2" = 2763;  "D2" = 131;  "E2" = 110;  "F2" = 816;  "G2" = 256
    "B3" = -0.01052456029596116; "C3" = -0.0321981526640378; "D3" = -0.05707477054479404; "E3" = -0.6687957480535223; "F3" = -0.001651116592686241; "G3" = 0.313326173972469
    "B4" = -0.03542559296103653; "C4" = -0.07016940602780486; "D4" = -3.294366255338225; "E4" = -15.53624877616672; "F4" = -0.1385574058898735; "G4" = -0.1152982524650981
    "B5" = 3.436401711423847; "C5" = 3.3612301705088; "D5" = 3.231524112799261; "E5" = 3.845895096191404; "F5" = 3.790450138822823; "G5" = 2.93017704243276
    "B6" = 74;  "C6" = 47;  "D6" = 2;  "E6" = 1;  "F6" = 23;  "G6" = 1
    "B7" = 0.4960832566928129; "C7" = 0.4341680724834786; "D7" = 0.5428263380630094; "E7" = 0.2316452809869561; "F7" = 0.6121035607998226; "G7" = 0.6846531514767573
    "B8" = 0.4836658966616074; "C8" = 0.4133529340087717; "D8" = -0.8572680016190244; "E8" = -6.613696761129253; "F8" = 0.5590856374502866; "G8" = 0.4878124434813574
    "B9" = 2.426665324080497; "C9" = 2.488630612230798; "D9" = 2.125177201954327; "E9" = 2.609616700909467; "F9" = 2.358796337740918; "G9" = 1.985696462501479
    "B10" = 72;  "C10" = 48;  "D10" = 3;  "E10" = 1;  "F10" = 19;  "G10" = 1
    "B11" = 0.5066078169887741; "C11" = 0.4663662251475164; "D11" = 0.5999011086078034; "E11" = 0.9004410290404784; "F11" = 0.6137546773925089; "G11" = 0.3713269775042883
    "B12" = 0.519091489622644; "C12" = 0.4835223400365766; "D12" = 2.437098253719201; "E12" = 8.92255201503747; "F12" = 0.6976430433401601; "G12" = 0.6031106959464555
    "B13" = -1.00973638734335; "C13" = -0.8725995582780022; "D13" = -1.106346910844934; "E13" = -1.236278395281937; "F13" = -1.431653801081905; "G13" = -0.9444805799312816
    "B14" = -2;  "C14" = 1;  "D14" = 1;  "E14" = 0;  "F14" = -4;  "G14" = 0
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}
